$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 41
$ws.Range("I41").Value = 41
$ws.Range("K41").Value = 41
$ws.Range("M41").Value = 399

$ws.Range("H45").Value = 4788.5713
$ws.Range("J45").Value = 6420.6
$ws.Range("L45").Value = 19261.8
$ws.Range("N45").Value = -19645.8

$ws.Range("H64").Value = 9911.25
$ws.Range("I64").Value = 9997
$ws.Range("J64").Value = 9894.1
$ws.Range("K64").Value = 9997
$ws.Range("L64").Value = 9894.1
$ws.Range("M64").Value = -9749
$ws.Range("N64").Value = -10390.1

$ws.Range("H67").Value = 9911.25
$ws.Range("I67").Value = 9997
$ws.Range("J67").Value = 9894.1
$ws.Range("K67").Value = 9997
$ws.Range("L67").Value = 9894.1
$ws.Range("M67").Value = -9139
$ws.Range("N67").Value = -11610.1

$ws.Range("H69").Value = 8144.5
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 8703.091
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 26109.273
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -27857.273

$ws.Range("H72").Value = 8144.5
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 8703.091
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 78327.819
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -87063.819

$ws.Range("H86").Value = 2923.5833
$ws.Range("J86").Value = 2956.8823
$ws.Range("L86").Value = 2956.8823
$ws.Range("N86").Value = -5202.8823

$ws.Range("H89").Value = 2923.5833
$ws.Range("J89").Value = 2956.8823
$ws.Range("L89").Value = 14784.4115
$ws.Range("N89").Value = -26016.4115

$ws.Range("H132").Value = 1564.0435
$ws.Range("I132").Value = 1557.8292
$ws.Range("K132").Value = 4673.487599999999
$ws.Range("M132").Value = -2143.487599999999

$ws.Range("H137").Value = 3686.8386
$ws.Range("I137").Value = 3184.4546
$ws.Range("K137").Value = 9553.363799999999
$ws.Range("M137").Value = -7003.363799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12271.6
$ws.Range("J2").Value = 16673.715
$ws.Range("L2").Value = 16673.715
$ws.Range("N2").Value = -16899.715

$ws.Range("H32").Value = 1780.9747
$ws.Range("I32").Value = 1643.7733
$ws.Range("K32").Value = 1643.7733
$ws.Range("M32").Value = -1356.7733

$ws.Range("H45").Value = 4650.933
$ws.Range("I45").Value = 2878.1428
$ws.Range("K45").Value = 2878.1428
$ws.Range("M45").Value = -2501.1428

$ws.Range("H61").Value = 5626.911
$ws.Range("I61").Value = 5609.8374
$ws.Range("K61").Value = 5609.8374
$ws.Range("M61").Value = -5397.8374

$ws.Range("H74").Value = 19610700
$ws.Range("I74").Value = 27780248
$ws.Range("K74").Value = 27780248
$ws.Range("M74").Value = -27779374

$ws.Range("H77").Value = 19610700
$ws.Range("I77").Value = 27780248
$ws.Range("K77").Value = 138901240
$ws.Range("M77").Value = -138896872

$ws.Range("H116").Value = 12271.6
$ws.Range("J116").Value = 16673.715
$ws.Range("L116").Value = 16673.715
$ws.Range("N116").Value = -21261.715

$ws.Range("H135").Value = 46657.637
$ws.Range("J135").Value = 46657.637
$ws.Range("L135").Value = 46657.637
$ws.Range("N135").Value = -56797.637

$ws.Range("H136").Value = 5626.911
$ws.Range("I136").Value = 5609.8374
$ws.Range("K136").Value = 16829.5122
$ws.Range("M136").Value = -14279.5122

$ws.Range("H137").Value = 69991.336
$ws.Range("J137").Value = 69991.336
$ws.Range("L137").Value = 69991.336
$ws.Range("N137").Value = -80191.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12271.6
$ws.Range("J3").Value = 16673.715
$ws.Range("L3").Value = 16673.715
$ws.Range("N3").Value = -16901.715

$ws.Range("H20").Value = 4510.6
$ws.Range("I20").Value = 3423.8333
$ws.Range("K20").Value = 3423.8333
$ws.Range("M20").Value = -3176.8333

$ws.Range("H113").Value = 5305.4443
$ws.Range("I113").Value = 5305.4443
$ws.Range("K113").Value = 5305.4443
$ws.Range("M113").Value = -3135.4443

$ws.Range("H134").Value = 1627.3334
$ws.Range("I134").Value = 1048.7273
$ws.Range("K134").Value = 3146.1819
$ws.Range("M134").Value = -611.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 42218.883
$ws.Range("I31").Value = 1283.15
$ws.Range("K31").Value = 1283.15
$ws.Range("M31").Value = -988.1500000000001

$ws.Range("H34").Value = 42218.883
$ws.Range("I34").Value = 1283.15
$ws.Range("K34").Value = 1283.15
$ws.Range("M34").Value = -1081.15

$ws.Range("H39").Value = 8566.666999999999
$ws.Range("I39").Value = 5350
$ws.Range("J39").Value = 15000
$ws.Range("K39").Value = 5350
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = -4959
$ws.Range("N39").Value = -15782

$ws.Range("H49").Value = 8566.666999999999
$ws.Range("I49").Value = 5350
$ws.Range("J49").Value = 15000
$ws.Range("K49").Value = 5350
$ws.Range("L49").Value = 15000
$ws.Range("M49").Value = -5168
$ws.Range("N49").Value = -15364

$ws.Range("H74").Value = 47999.3
$ws.Range("J74").Value = 65998.60000000001
$ws.Range("L74").Value = 65998.60000000001
$ws.Range("N74").Value = -67746.60000000001

$ws.Range("H77").Value = 47999.3
$ws.Range("J77").Value = 65998.60000000001
$ws.Range("L77").Value = 197995.8
$ws.Range("N77").Value = -206731.8

$ws.Range("H122").Value = 47703.832
$ws.Range("I122").Value = 61876.47
$ws.Range("K122").Value = 185629.41
$ws.Range("M122").Value = -183179.41

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4632311
$ws.Range("I129").Value = 650.8
$ws.Range("J129").Value = 10421886
$ws.Range("K129").Value = 1952.4
$ws.Range("L129").Value = 31265658
$ws.Range("M129").Value = 3047.6
$ws.Range("N129").Value = -31275658

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 55000
$ws.Range("J80").Value = 55000
$ws.Range("L80").Value = 55000
$ws.Range("N80").Value = -57246

$ws.Range("H83").Value = 55000
$ws.Range("J83").Value = 55000
$ws.Range("L83").Value = 165000
$ws.Range("N83").Value = -176232

$ws.Range("H136").Value = 11876.883
$ws.Range("I136").Value = 5322.5557
$ws.Range("J136").Value = 19250.5
$ws.Range("K136").Value = 15967.6671
$ws.Range("L136").Value = 57751.5
$ws.Range("M136").Value = -13417.6671
$ws.Range("N136").Value = -62851.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2027.5319
$ws.Range("I136").Value = 1227.6511
$ws.Range("K136").Value = 3682.9533
$ws.Range("M136").Value = -1132.9533
